$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Resolving-Mac" target-cluster rows (previously rows 5 and 9)
$ws.Rows("9").Delete()
$ws.Rows("5").Delete()

# Update cell values for rows 2-7 with the refreshed TPM-based statistics
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.601971
$ws.Range("H2").Value = 1.805913
$ws.Range("I2").Value = 0.09594307528308157
$ws.Range("J2").Value = 0.09594307528308157
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.261293
$ws.Range("N2").Value = 0.783879
$ws.Range("O2").Value = 0.0361900776512412
$ws.Range("P2").Value = 0.03619007765124121
$ws.Range("Q2").Value = 0.157290808503
$ws.Range("R2").Value = 1.415617276527
$ws.Range("S2").Value = 0.003472187344593603
$ws.Range("T2").Value = 0.003472187344593603

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.601971
$ws.Range("H3").Value = 1.805913
$ws.Range("I3").Value = 0.09594307528308157
$ws.Range("J3").Value = 0.09594307528308157
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7775033333333333
$ws.Range("N3").Value = 2.33251
$ws.Range("O3").Value = 0.1076871787894517
$ws.Range("P3").Value = 0.1076871787894517
$ws.Range("Q3").Value = 0.46803445907
$ws.Range("R3").Value = 4.212310131630001
$ws.Range("S3").Value = 0.01033183910161903
$ws.Range("T3").Value = 0.01033183910161903

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.601971
$ws.Range("H4").Value = 1.805913
$ws.Range("I4").Value = 0.09594307528308157
$ws.Range("J4").Value = 0.09594307528308157
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.181221333333333
$ws.Range("N4").Value = 18.543664
$ws.Range("O4").Value = 0.8561227435593071
$ws.Range("P4").Value = 0.8561227435593072
$ws.Range("Q4").Value = 3.720915987248
$ws.Range("R4").Value = 33.488243885232
$ws.Range("S4").Value = 0.08213904883686894
$ws.Range("T4").Value = 0.08213904883686896

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.672280666666667
$ws.Range("H5").Value = 17.016842
$ws.Range("I5").Value = 0.9040569247169185
$ws.Range("J5").Value = 0.9040569247169185
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.261293
$ws.Range("N5").Value = 0.783879
$ws.Range("O5").Value = 0.0361900776512412
$ws.Range("P5").Value = 0.03619007765124121
$ws.Range("Q5").Value = 1.482127232235333
$ws.Range("R5").Value = 13.339145090118
$ws.Range("S5").Value = 0.0327178903066476
$ws.Range("T5").Value = 0.0327178903066476

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.672280666666667
$ws.Range("H6").Value = 17.016842
$ws.Range("I6").Value = 0.9040569247169185
$ws.Range("J6").Value = 0.9040569247169185
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7775033333333333
$ws.Range("N6").Value = 2.33251
$ws.Range("O6").Value = 0.1076871787894517
$ws.Range("P6").Value = 0.1076871787894517
$ws.Range("Q6").Value = 4.410217125935556
$ws.Range("R6").Value = 39.69195413342
$ws.Range("S6").Value = 0.09735533968783268
$ws.Range("T6").Value = 0.0973553396878327

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.672280666666667
$ws.Range("H7").Value = 17.016842
$ws.Range("I7").Value = 0.9040569247169185
$ws.Range("J7").Value = 0.9040569247169185
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.181221333333333
$ws.Range("N7").Value = 18.543664
$ws.Range("O7").Value = 0.8561227435593071
$ws.Range("P7").Value = 0.8561227435593072
$ws.Range("Q7").Value = 35.06162226545422
$ws.Range("R7").Value = 315.554600389088
$ws.Range("S7").Value = 0.7739836947224381
$ws.Range("T7").Value = 0.7739836947224382
